$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates to match the scraped cryptocurrency data refresh.
$ws.Range("D2").Value = "66.866.94"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").Value = "3.100.90"
$ws.Range("E3").Value = "  +5.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.96"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.58"
$ws.Range("E6").Value = "  +5.63%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.097.23"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("E10").Value = "  +3.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.42"
$ws.Range("E11").Value = "  -3.93%  "
$ws.Range("E12").Value = "  +4.31%  "
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.22"
$ws.Range("E14").Value = "  +6.68%  "
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "3.612.47"
$ws.Range("E16").Value = "  +5.18%  "
$ws.Range("D17").Value = "66.866.01"
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").Value = "3.102.59"
$ws.Range("E19").Value = "  +5.26%  "
$ws.Range("E20").Value = "  +2.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "479.94"
$ws.Range("E21").Value = "  +7.57%  "
$ws.Range("E22").Value = "  +2.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.91"
$ws.Range("E24").Value = "  +2.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.20"
$ws.Range("E25").Value = "  +8.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.37"
$ws.Range("E26").Value = "  +4.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.05"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("E31").Value = "  +3.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.85"
$ws.Range("E32").Value = "  +6.01%  "
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.89"
$ws.Range("E36").Value = "  +3.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.994"
$ws.Range("E37").Value = "  +2.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.05"
$ws.Range("E38").Value = "  +2.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.11"
$ws.Range("E39").Value = "  +6.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.14"
$ws.Range("E40").Value = "  +1.95%  "
$ws.Range("E41").Value = "  +4.19%  "
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.838.74"
$ws.Range("E45").Value = "  +6.15%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0361"
$ws.Range("E46").Value = "  +2.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "386.01"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.95"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.89"
$ws.Range("E50").Value = "  +3.74%  "
$ws.Range("E51").Value = "  +2.80%  "
